$wb = $excel.ActiveWorkbook

# Sheet "展览" - update 想去人数 (want-to-go count) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 5596
$wsExhibit.Range("F9").Value = 4

# Sheet "全部类型" - same underlying rows, update corresponding values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5596
$wsAll.Range("F10").Value = 4
